{"js": "const replacements = [\n  [\"916\u00d73=2748\", \"238\u00d79=2142\"],\n  [\"892\u00d72=1784\", \"885\u00d76=5310\"],\n  [\"741\u00d72=1482\", \"589\u00d76=3534\"],\n  [\"830\u00d79=7470\", \"734\u00d73=2202\"],\n  [\"294\u00d73=882\", \"685\u00d78=5480\"],\n  [\"333\u00d78=2664\", \"592\u00d79=5328\"],\n  [\"547\u00d72=1094\", \"476\u00d79=4284\"],\n  [\"291\u00d77=2037\", \"923\u00d77=6461\"],\n  [\"858\u00d72=1716\", \"843\u00d73=2529\"],\n  [\"855\u00d73=2565\", \"827\u00d74=3308\"],\n  [\"485\u00d74=1940\", \"486\u00d72=972\"],\n  [\"318\u00d77=2226\", \"785\u00d72=1570\"],\n  [\"781\u00d76=4686\", \"699\u00d72=1398\"],\n  [\"224\u00d72=448\", \"736\u00d73=2208\"],\n  [\"364\u00d75=1820\", \"723\u00d78=5784\"],\n  [\"140\u00d79=1260\", \"354\u00d77=2478\"],\n  [\"464\u00d74=1856\", \"524\u00d74=2096\"],\n  [\"308\u00d76=1848\", \"761\u00d76=4566\"],\n  [\"546\u00d79=4914\", \"407\u00d73=1221\"],\n  [\"348\u00d78=2784\", \"366\u00d79=3294\"],\n  [\"111\u00d77=777\", \"125\u00d78=1000\"],\n  [\"611\u00d75=3055\", \"825\u00d79=7425\"],\n  [\"128\u00d75=640\", \"146\u00d74=584\"],\n  [\"382\u00d79=3438\", \"330\u00d72=660\"],\n  [\"601\u00d78=4808\", \"152\u00d72=304\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"916\u00d73=2748\", \"238\u00d79=2142\"),\n  @(\"892\u00d72=1784\", \"885\u00d76=5310\"),\n  @(\"741\u00d72=1482\", \"589\u00d76=3534\"),\n  @(\"830\u00d79=7470\", \"734\u00d73=2202\"),\n  @(\"294\u00d73=882\", \"685\u00d78=5480\"),\n  @(\"333\u00d78=2664\", \"592\u00d79=5328\"),\n  @(\"547\u00d72=1094\", \"476\u00d79=4284\"),\n  @(\"291\u00d77=2037\", \"923\u00d77=6461\"),\n  @(\"858\u00d72=1716\", \"843\u00d73=2529\"),\n  @(\"855\u00d73=2565\", \"827\u00d74=3308\"),\n  @(\"485\u00d74=1940\", \"486\u00d72=972\"),\n  @(\"318\u00d77=2226\", \"785\u00d72=1570\"),\n  @(\"781\u00d76=4686\", \"699\u00d72=1398\"),\n  @(\"224\u00d72=448\", \"736\u00d73=2208\"),\n  @(\"364\u00d75=1820\", \"723\u00d78=5784\"),\n  @(\"140\u00d79=1260\", \"354\u00d77=2478\"),\n  @(\"464\u00d74=1856\", \"524\u00d74=2096\"),\n  @(\"308\u00d76=1848\", \"761\u00d76=4566\"),\n  @(\"546\u00d79=4914\", \"407\u00d73=1221\"),\n  @(\"348\u00d78=2784\", \"366\u00d79=3294\"),\n  @(\"111\u00d77=777\", \"125\u00d78=1000\"),\n  @(\"611\u00d75=3055\", \"825\u00d79=7425\"),\n  @(\"128\u00d75=640\", \"146\u00d74=584\"),\n  @(\"382\u00d79=3438\", \"330\u00d72=660\"),\n  @(\"601\u00d78=4808\", \"152\u00d72=304\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
